$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing study session row (row 2) with refreshed timer values
$ws.Range("C2").Value = 45688.88619232639
$ws.Range("D2").Value = 45688.88624818287

# Append a new study session row (row 3) recorded by the Studying timer
$ws.Range("A3").Value = "MAT141"
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 45688.88619232506
$ws.Range("D3").Value = 45688.88624818008

# Match the date/time number format already used by the start/end columns
$ws.Range("C3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
